$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue number + week date range) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Crime statistics table updates (rows 16-29) ---
$ws.Range("I14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -86.666666666666
$ws.Range("I16").Value = 120
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = 4.347826086956
$ws.Range("L16").Value = 27.659574468085
$ws.Range("M16").Value = -20.529801324503
$ws.Range("N16").Value = -89.218328840970
$ws.Range("C17").Value = 3
$ws.Range("C14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "***.*"
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 283.333333333333
$ws.Range("I17").Value = 154
$ws.Range("K17").Value = 11.594202898550
$ws.Range("L17").Value = 7.692307692307
$ws.Range("M17").Value = 62.105263157894
$ws.Range("N17").Value = -59.473684210526
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -31.578947368421
$ws.Range("I18").Value = 163
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = 32.520325203252
$ws.Range("L18").Value = 1.242236024844
$ws.Range("M18").Value = 61.386138613861
$ws.Range("N18").Value = -76.879432624113
$ws.Range("C19").Value = 14
$ws.Range("E19").Value = 7.692307692307
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 16.666666666666
$ws.Range("I19").Value = 547
$ws.Range("J19").Value = 439
$ws.Range("K19").Value = 24.601366742596
$ws.Range("L19").Value = 58.550724637681
$ws.Range("M19").Value = 26.620370370370
$ws.Range("N19").Value = -42.238648363252
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").Value = 8
$ws.Range("H20").Value = 60
$ws.Range("L20").Value = 34.146341463414
$ws.Range("N20").Value = -91.352201257861
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -8.695652173913
$ws.Range("F21").Value = 96
$ws.Range("H21").Value = 10.344827586206
$ws.Range("I21").Value = 1050
$ws.Range("J21").Value = 859
$ws.Range("K21").Value = 22.235157159487
$ws.Range("L21").Value = 32.743362831858
$ws.Range("M21").Value = 26.201923076923
$ws.Range("N21").Value = -72.440944881889
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 53
$ws.Range("J22").Value = 44
$ws.Range("K22").Value = 20.454545454545
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -14.516129032258
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 32
$ws.Range("J23").Value = 41
$ws.Range("K23").Value = -21.951219512195
$ws.Range("L23").Value = 6.666666666666
$ws.Range("M23").Value = 39.130434782608
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = 42.424242424242
$ws.Range("F24").Value = 180
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = 45.161290322580
$ws.Range("I24").Value = 1666
$ws.Range("J24").Value = 1141
$ws.Range("K24").Value = 46.012269938650
$ws.Range("L24").Value = 105.679012345679
$ws.Range("M24").Value = 27.467482785003
$ws.Range("C25").Value = 14
$ws.Range("E25").Value = 55.555555555555
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 25.806451612903
$ws.Range("I25").Value = 320
$ws.Range("J25").Value = 229
$ws.Range("K25").Value = 39.737991266375
$ws.Range("L25").Value = 33.891213389121
$ws.Range("M25").Value = 4.918032786885
$ws.Range("F26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 59
$ws.Range("J27").Value = 58
$ws.Range("K27").Value = 1.724137931034
$ws.Range("L27").Value = 59.459459459459
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = "0"
$ws.Range("F28").Value = 1
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = "0"
$ws.Range("F29").Value = 1
